$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Corp Trade Compliance")
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0.004
$ws.Range("P4").Value = 0.004
$ws.Range("Q4").Value = 0.012
$ws.Range("R4").Value = 0.004
$ws.Range("S4").Value = 0.004
$ws.Range("T4").Value = 0.004
$ws.Range("U4").Value = 0.012
$ws.Range("V4").Value = 0.048

$ws = $wb.Worksheets.Item("Finance")
$ws.Range("D2").Value = 0.0613
$ws.Range("D3").Value = 0.0613
$ws.Range("D4").Value = 0.0613
$ws.Range("F4").Value = 0.0071
$ws.Range("I4").Value = 0.019
$ws.Range("J4").Value = 0.0266
$ws.Range("L4").Value = 0.0046
$ws.Range("M4").Value = 0.0331
$ws.Range("N4").Value = 0.0092
$ws.Range("O4").Value = 0.00875833333333333
$ws.Range("P4").Value = 0.00875833333333333
$ws.Range("Q4").Value = 0.026275
$ws.Range("R4").Value = 0.00875833333333333
$ws.Range("S4").Value = 0.00875833333333333
$ws.Range("T4").Value = 0.00875833333333333
$ws.Range("U4").Value = 0.026275
$ws.Range("V4").Value = 0.1051
$ws.Range("D5").Value = 0.333333333333333
$ws.Range("D6").Value = 0.333333333333333
$ws.Range("D7").Value = 0.333333333333333
$ws.Range("L7").Value = 0.4
$ws.Range("M7").Value = 0.25
$ws.Range("N7").Value = 1
$ws.Range("O7").Value = 0.333333333333333
$ws.Range("P7").Value = 0.333333333333333
$ws.Range("Q7").Value = 0.333333333333333
$ws.Range("R7").Value = 0.333333333333333
$ws.Range("S7").Value = 0.333333333333333
$ws.Range("T7").Value = 0.333333333333333
$ws.Range("U7").Value = 0.333333333333333
$ws.Range("V7").Value = 0.333333333333333

$ws = $wb.Worksheets.Item("L1_Corporate")
$ws.Range("D2").Value = 0.0683
$ws.Range("D3").Value = 0.0683
$ws.Range("D4").Value = 0.0683
$ws.Range("F4").Value = 0.0051
$ws.Range("I4").Value = 0.0321
$ws.Range("M4").Value = 0.0291
$ws.Range("N4").Value = 0.0072
$ws.Range("O4").Value = 0.00975833333333333
$ws.Range("P4").Value = 0.00975833333333333
$ws.Range("Q4").Value = 0.029275
$ws.Range("R4").Value = 0.00975833333333333
$ws.Range("S4").Value = 0.00975833333333333
$ws.Range("T4").Value = 0.00975833333333333
$ws.Range("U4").Value = 0.029275
$ws.Range("V4").Value = 0.1171
$ws.Range("D5").Value = 0.444444444444444
$ws.Range("D6").Value = 0.444444444444444
$ws.Range("D7").Value = 0.444444444444444
$ws.Range("L7").Value = 0.5714
$ws.Range("M7").Value = 0.4074
$ws.Range("N7").Value = 0.6667
$ws.Range("O7").Value = 0.444444444444444
$ws.Range("P7").Value = 0.444444444444444
$ws.Range("Q7").Value = 0.444444444444444
$ws.Range("R7").Value = 0.444444444444444
$ws.Range("S7").Value = 0.444444444444444
$ws.Range("T7").Value = 0.444444444444444
$ws.Range("U7").Value = 0.444444444444444
$ws.Range("V7").Value = 0.444444444444444

$ws = $wb.Worksheets.Item("Corp Business Development")
$ws.Range("N7").Value = $null

$ws = $wb.Worksheets.Item("Corp EH&S")
$ws.Range("D2").Value = 0.339
$ws.Range("D3").Value = 0.339
$ws.Range("D4").Value = 0.339
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0.048425
$ws.Range("P4").Value = 0.048425
$ws.Range("Q4").Value = 0.145275
$ws.Range("R4").Value = 0.048425
$ws.Range("S4").Value = 0.048425
$ws.Range("T4").Value = 0.048425
$ws.Range("U4").Value = 0.145275
$ws.Range("V4").Value = 0.5811
$ws.Range("N5").Value = $null

$ws = $wb.Worksheets.Item("Corp Information Technology")
$ws.Range("D2").Value = 0.0523
$ws.Range("D3").Value = 0.0523
$ws.Range("D4").Value = 0.0523
$ws.Range("H4").Value = 0.0131
$ws.Range("I4").Value = 0.0291
$ws.Range("J4").Value = 0.0098
$ws.Range("M4").Value = 0.0231
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0.007475
$ws.Range("P4").Value = 0.007475
$ws.Range("Q4").Value = 0.022425
$ws.Range("R4").Value = 0.007475
$ws.Range("S4").Value = 0.007475
$ws.Range("T4").Value = 0.007475
$ws.Range("U4").Value = 0.022425
$ws.Range("V4").Value = 0.0897
$ws.Range("N7").Value = $null

$ws = $wb.Worksheets.Item("Corp Legal")
$ws.Range("D2").Value = 0.2251
$ws.Range("D3").Value = 0.2251
$ws.Range("D4").Value = 0.2251
$ws.Range("G4").Value = 0.0606
$ws.Range("H4").Value = 0.0667
$ws.Range("I4").Value = 0.1223
$ws.Range("J4").Value = 0.069
$ws.Range("M4").Value = 0.0673
$ws.Range("N4").Value = 0.0323
$ws.Range("O4").Value = 0.0321583333333333
$ws.Range("P4").Value = 0.0321583333333333
$ws.Range("Q4").Value = 0.096475
$ws.Range("R4").Value = 0.0321583333333333
$ws.Range("S4").Value = 0.0321583333333333
$ws.Range("T4").Value = 0.0321583333333333
$ws.Range("U4").Value = 0.096475
$ws.Range("V4").Value = 0.3859
$ws.Range("D5").Value = 0.166666666666667
$ws.Range("D6").Value = 0.166666666666667
$ws.Range("D7").Value = 0.166666666666667
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 0.166666666666667
$ws.Range("P7").Value = 0.166666666666667
$ws.Range("Q7").Value = 0.166666666666667
$ws.Range("R7").Value = 0.166666666666667
$ws.Range("S7").Value = 0.166666666666667
$ws.Range("T7").Value = 0.166666666666667
$ws.Range("U7").Value = 0.166666666666667
$ws.Range("V7").Value = 0.166666666666667

$ws = $wb.Worksheets.Item("Corp Logistics")
$ws.Range("D2").Value = 0.0803
$ws.Range("D3").Value = 0.0803
$ws.Range("D4").Value = 0.0803
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0.011475
$ws.Range("P4").Value = 0.011475
$ws.Range("Q4").Value = 0.034425
$ws.Range("R4").Value = 0.011475
$ws.Range("S4").Value = 0.011475
$ws.Range("T4").Value = 0.011475
$ws.Range("U4").Value = 0.034425
$ws.Range("V4").Value = 0.1377
$ws.Range("N7").Value = $null

$ws = $wb.Worksheets.Item("Corp RBS")
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0.0865833333333333
$ws.Range("P4").Value = 0.0865833333333333
$ws.Range("Q4").Value = 0.25975
$ws.Range("R4").Value = 0.0865833333333333
$ws.Range("S4").Value = 0.0865833333333333
$ws.Range("T4").Value = 0.0865833333333333
$ws.Range("U4").Value = 0.25975
$ws.Range("V4").Value = 1.039
$ws.Range("N7").Value = $null

$ws = $wb.Worksheets.Item("Corp Sourcing")
$ws.Range("D2").Value = 0.0946
$ws.Range("D3").Value = 0.0946
$ws.Range("D4").Value = 0.0946
$ws.Range("F4").Value = 0.0137
$ws.Range("H4").Value = 0.0429
$ws.Range("I4").Value = 0.0556
$ws.Range("J4").Value = 0.0137
$ws.Range("M4").Value = 0.0133
$ws.Range("N4").Value = 0.026
$ws.Range("O4").Value = 0.0135166666666667
$ws.Range("P4").Value = 0.0135166666666667
$ws.Range("Q4").Value = 0.04055
$ws.Range("R4").Value = 0.0135166666666667
$ws.Range("S4").Value = 0.0135166666666667
$ws.Range("T4").Value = 0.0135166666666667
$ws.Range("U4").Value = 0.04055
$ws.Range("V4").Value = 0.1622
$ws.Range("N7").Value = $null
